$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the libraryDate value for rows 13-17 (A13:A17) from "01.06.20" to "01.16.20"
# Force text number format first so Excel doesn't auto-convert the date-like string
# into a serial date, then clear the formatting back to the default (General) so the
# cells keep looking exactly as they did before (no explicit style index).
$ws.Range("A13:A17").NumberFormat = "@"
$ws.Range("A13:A17").Value = "01.16.20"
$ws.Range("A13:A17").ClearFormats()

# Update the active cell / selection on the sheet to A18 (single cell, no range)
$ws.Range("A18").Select()
